# Add a new forecast column (published 2020-05-12) and a new observed-date
# row (2020-05-26) to both the "cases" and "deaths" sheets, plus fill in
# the previously-missing "Observed" value for the 2020-05-12 row (B31).

$wb = $excel.ActiveWorkbook

# AG-column (col 33) diagonal values, keyed by row number, per sheet.
$casesAG = @{
    32 = 19426
    33 = 20247
    34 = 20765
    35 = 21518
    36 = 22046
    37 = 22456
    38 = 22964
    39 = 23444
    40 = 24146
    41 = 24605
    42 = 25165
    43 = 25599
    44 = 26000
}
$casesAG45 = 26583
$casesB31  = 18486

$deathsAG = @{
    32 = 2076
    33 = 2225
    34 = 2326
    35 = 2456
    36 = 2549
    37 = 2642
    38 = 2727
    39 = 2832
    40 = 2948
    41 = 3063
    42 = 3181
    43 = 3270
    44 = 3347
}
$deathsAG45 = 3463
$deathsB31  = 1928

$sheetNames = @("cases", "deaths")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "cases") {
        $agValues = $casesAG
        $ag45 = $casesAG45
        $b31 = $casesB31
    } else {
        $agValues = $deathsAG
        $ag45 = $deathsAG45
        $b31 = $deathsB31
    }

    # --- New column AG: header is the forecast-publish date 2020-05-12.
    # Force text so Excel doesn't coerce the date-looking string into a
    # date serial number, then clear the temporary formatting so the cell
    # keeps the default (unstyled) look.
    $ws.Range("AG1").NumberFormat = "@"
    $ws.Range("AG1").Value = "2020-05-12"
    $ws.Range("AG1").ClearFormats()

    # Rows 2-31: AG stays blank (no forecast yet for these older rows).
    # Rows 32-44: AG gets the diagonal forecast value.
    foreach ($row in $agValues.Keys) {
        $ws.Cells.Item($row, 33).Value = $agValues[$row]
    }

    # --- Previously missing "Observed" figure for the 2020-05-12 row.
    $ws.Range("B31").Value = $b31

    # --- New row 45: observed date 2020-05-26, with only the AG (2020-05-12
    # forecast) figure populated.
    $ws.Range("A45").NumberFormat = "@"
    $ws.Range("A45").Value = "2020-05-26"
    $ws.Range("A45").ClearFormats()

    $ws.Cells.Item(45, 33).Value = $ag45
}
